$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.19667366061813
$ws.Range("C2").Value = 8.35966439233923
$ws.Range("D2").Value = 9.420009377179049
$ws.Range("F2").Value = 35.59575027298209
$ws.Range("G2").Value = 39.04177226840294
$ws.Range("H2").Value = 16.5200001671281
$ws.Range("J2").Value = 10.92683161462506
$ws.Range("M2").Value = 18.16309274868948
$ws.Range("N2").Value = 18.29741856072221

$ws.Range("B3").Value = 14.71976926512402
$ws.Range("C3").Value = 7.916798218922285
$ws.Range("D3").Value = 9.417146225209748
$ws.Range("F3").Value = 35.53489286404012
$ws.Range("G3").Value = 38.81737476030996
$ws.Range("H3").Value = 16.54545636646835
$ws.Range("J3").Value = 10.95170149816921
$ws.Range("M3").Value = 18.02749359007558
$ws.Range("N3").Value = 18.36915314031053

$ws.Range("B4").Value = 14.42253728063235
$ws.Range("C4").Value = 7.633067536066225
$ws.Range("D4").Value = 9.4164360498602
$ws.Range("F4").Value = 35.50791891676248
$ws.Range("G4").Value = 38.69356945837765
$ws.Range("H4").Value = 16.56488238752863
$ws.Range("J4").Value = 10.96839731360714
$ws.Range("M4").Value = 17.94749226624987
$ws.Range("N4").Value = 18.41513602071812

$ws.Range("B5").Value = 14.30050635823942
$ws.Range("C5").Value = 7.514607895658941
$ws.Range("D5").Value = 9.41641116639007
$ws.Range("F5").Value = 35.4995436465108
$ws.Range("G5").Value = 38.64667154836192
$ws.Range("H5").Value = 16.5737504877862
$ws.Range("J5").Value = 10.97555951393209
$ws.Range("M5").Value = 17.91573913644444
$ws.Range("N5").Value = 18.43436335972076

$ws.Range("B6").Value = 14.28019466399305
$ws.Range("C6").Value = 7.494770754708786
$ws.Range("D6").Value = 9.416423039672758
$ws.Range("F6").Value = 35.49831103854073
$ws.Range("G6").Value = 38.63909984055984
$ws.Range("H6").Value = 16.57528043472299
$ws.Range("J6").Value = 10.97677044663457
$ws.Range("M6").Value = 17.91051859343004
$ws.Range("N6").Value = 18.43758562631598

$ws.Range("B7").Value = 14.42089492116144
$ws.Range("C7").Value = 7.631481242550177
$ws.Range("D7").Value = 9.416434641924416
$ws.Range("F7").Value = 35.50779536698776
$ws.Range("G7").Value = 38.69292254090673
$ws.Range("H7").Value = 16.56499813564169
$ws.Range("J7").Value = 10.96849245386837
$ws.Range("M7").Value = 17.94706056144567
$ws.Range("N7").Value = 18.41539334537183

$ws.Range("B8").Value = 15.03326964708102
$ws.Range("C8").Value = 8.209487385395796
$ws.Range("D8").Value = 9.418805331450981
$ws.Range("F8").Value = 35.57261252028404
$ws.Range("G8").Value = 38.96152513993984
$ws.Range("H8").Value = 16.52798813394136
$ws.Range("J8").Value = 10.9351108914562
$ws.Range("M8").Value = 18.11568000214562
$ws.Range("N8").Value = 18.32175154854459

$ws.Range("B9").Value = 16.19097119276024
$ws.Range("C9").Value = 9.244549669876264
$ws.Range("D9").Value = 9.431720642148997
$ws.Range("F9").Value = 35.7819076529858
$ws.Range("G9").Value = 39.59708061383115
$ws.Range("H9").Value = 16.48563705128067
$ws.Range("J9").Value = 10.88095964468478
$ws.Range("M9").Value = 18.47084485801577
$ws.Range("N9").Value = 18.15341369537461

$ws.Range("B10").Value = 17.0053498906668
$ws.Range("C10").Value = 9.940016236098085
$ws.Range("D10").Value = 9.446182342305946
$ws.Range("F10").Value = 35.9852423891647
$ws.Range("G10").Value = 40.12715428272708
$ws.Range("H10").Value = 16.47307811007994
$ws.Range("J10").Value = 10.84806743170252
$ws.Range("M10").Value = 18.74484426842238
$ws.Range("N10").Value = 18.03894480612465

$ws.Range("B11").Value = 17.36619142497605
$ws.Range("C11").Value = 10.24153779162444
$ws.Range("D11").Value = 9.453825441684263
$ws.Range("F11").Value = 36.08833566326475
$ws.Range("G11").Value = 40.38119980267168
$ws.Range("H11").Value = 16.471415308326
$ws.Range("J11").Value = 10.83460008807804
$ws.Range("M11").Value = 18.87191032456503
$ws.Range("N11").Value = 17.9888450584778

$ws.Range("B12").Value = 17.50132014298477
$ws.Range("C12").Value = 10.35353373150059
$ws.Range("D12").Value = 9.456871254017685
$ws.Range("F12").Value = 36.12887915967967
$ws.Range("G12").Value = 40.47918162166923
$ws.Range("H12").Value = 16.47136919630642
$ws.Range("J12").Value = 10.82971536304241
$ws.Range("M12").Value = 18.92033904623554
$ws.Range("N12").Value = 17.97015543184076

$ws.Range("B13").Value = 17.4722870245294
$ws.Range("C13").Value = 10.32951119156154
$ws.Range("D13").Value = 9.456208569311912
$ws.Range("F13").Value = 36.12008082220719
$ws.Range("G13").Value = 40.45800156825986
$ws.Range("H13").Value = 16.47135315959082
$ws.Range("J13").Value = 10.83075781070503
$ws.Range("M13").Value = 18.90989578694044
$ws.Range("N13").Value = 17.97416805823238

$ws.Range("B14").Value = 17.3773395884078
$ws.Range("C14").Value = 10.25079577276911
$ws.Range("D14").Value = 9.454072994488136
$ws.Range("F14").Value = 36.09164117214301
$ws.Range("G14").Value = 40.38922559958945
$ws.Range("H14").Value = 16.47139981422291
$ws.Range("J14").Value = 10.83419390851027
$ws.Range("M14").Value = 18.87588855235876
$ws.Range("N14").Value = 17.98730180795408

$ws.Range("B15").Value = 17.31898072755625
$ws.Range("C15").Value = 10.20229465098297
$ws.Range("D15").Value = 9.452784581204551
$ws.Range("F15").Value = 36.07441635442339
$ws.Range("G15").Value = 40.3473278356195
$ws.Range("H15").Value = 16.47150441469599
$ws.Range("J15").Value = 10.83632662587031
$ws.Range("M15").Value = 18.85509761635944
$ws.Range("N15").Value = 17.99538329370983

$ws.Range("B16").Value = 16.98156294039891
$ws.Range("C16").Value = 9.920007918230267
$ws.Range("D16").Value = 9.445704137881494
$ws.Range("F16").Value = 35.97871661724837
$ws.Range("G16").Value = 40.11080520298466
$ws.Range("H16").Value = 16.47326838530094
$ws.Range("J16").Value = 10.84897764892281
$ws.Range("M16").Value = 18.73658581309251
$ws.Range("N16").Value = 18.04225850046291

$ws.Range("B17").Value = 16.77200607099881
$ws.Range("C17").Value = 9.742993163709999
$ws.Range("D17").Value = 9.441632075496166
$ws.Range("F17").Value = 35.92270875221491
$ws.Range("G17").Value = 39.96896040164992
$ws.Range("H17").Value = 16.47538882383608
$ws.Range("J17").Value = 10.85712168660071
$ws.Range("M17").Value = 18.6644770599682
$ws.Range("N17").Value = 18.07151903160343

$ws.Range("B18").Value = 16.65057557147543
$ws.Range("C18").Value = 9.63978343197879
$ws.Range("D18").Value = 9.439390215478403
$ws.Range("F18").Value = 35.89149311116563
$ws.Range("G18").Value = 39.88859527656238
$ws.Range("H18").Value = 16.47698963546441
$ws.Range("J18").Value = 10.86194668377133
$ws.Range("M18").Value = 18.62323224842499
$ws.Range("N18").Value = 18.08853469337453

$ws.Range("B19").Value = 16.60931109982331
$ws.Range("C19").Value = 9.604600345334195
$ws.Range("D19").Value = 9.438648428901633
$ws.Range("F19").Value = 35.88109608654859
$ws.Range("G19").Value = 39.86159686278594
$ws.Range("H19").Value = 16.47759707253159
$ws.Range("J19").Value = 10.86360452028751
$ws.Range("M19").Value = 18.60930812052117
$ws.Range("N19").Value = 18.0943278609795

$ws.Range("B20").Value = 16.79440778186057
$ws.Range("C20").Value = 9.761981499817789
$ws.Range("D20").Value = 9.442055185402113
$ws.Range("F20").Value = 35.92856766519151
$ws.Range("G20").Value = 39.9839342716931
$ws.Range("H20").Value = 16.4751236384989
$ws.Range("J20").Value = 10.85624017142775
$ws.Range("M20").Value = 18.67212959639937
$ws.Range("N20").Value = 18.06838498258659

$ws.Range("B21").Value = 17.40527001467636
$ws.Range("C21").Value = 10.27397600253451
$ws.Range("D21").Value = 9.454696164458642
$ws.Range("F21").Value = 36.09995392669598
$ws.Range("G21").Value = 40.40937908301471
$ws.Range("H21").Value = 16.47137026597358
$ws.Range("J21").Value = 10.83317880638411
$ws.Range("M21").Value = 18.88586914060608
$ws.Range("N21").Value = 17.98343646507712

$ws.Range("B22").Value = 17.79562745072377
$ws.Range("C22").Value = 10.59584964567958
$ws.Range("D22").Value = 9.463840282820998
$ws.Range("F22").Value = 36.22072227742987
$ws.Range("G22").Value = 40.69777038936682
$ws.Range("H22").Value = 16.47231888591872
$ws.Range("J22").Value = 10.81936044015566
$ws.Range("M22").Value = 19.02735904586238
$ws.Range("N22").Value = 17.92956101948286

$ws.Range("B23").Value = 17.58813801269841
$ws.Range("C23").Value = 10.42523927712529
$ws.Range("D23").Value = 9.458879668973454
$ws.Range("F23").Value = 36.15547165712722
$ws.Range("G23").Value = 40.54293089877603
$ws.Range("H23").Value = 16.47150106084229
$ws.Range("J23").Value = 10.82662085751509
$ws.Range("M23").Value = 18.95169086760206
$ws.Range("N23").Value = 17.95816552603734

$ws.Range("B24").Value = 16.78428292815181
$ws.Range("C24").Value = 9.753401359565537
$ws.Range("D24").Value = 9.441863588242088
$ws.Range("F24").Value = 35.92591578405492
$ws.Range("G24").Value = 39.97716089331391
$ws.Range("H24").Value = 16.47524233976932
$ws.Range("J24").Value = 10.85663825938192
$ws.Range("M24").Value = 18.66866922318114
$ws.Range("N24").Value = 18.0698012836081

$ws.Range("B25").Value = 15.88348159461778
$ws.Range("C25").Value = 8.975652734474664
$ws.Range("D25").Value = 9.427348291375788
$ws.Range("F25").Value = 35.71653362335272
$ws.Range("G25").Value = 39.41380296611813
$ws.Range("H25").Value = 16.49384356895927
$ws.Range("J25").Value = 10.89439838148963
$ws.Range("M25").Value = 18.37232993137215
$ws.Range("N25").Value = 18.19732813272803
